# Fruta / hortaliza, semanal
# Insert a new weekly record at row 40 (pushing the existing rows 40-124
# down to 41-125) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40:124 down by one row, creating a blank row 40.
$ws.Range("A40:R40").Insert()

# Populate the newly inserted row 40 with the new record.
$ws.Range("A40").Value = 2
$ws.Range("B40").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 45259
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 7000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = 7500
$ws.Range("N40").Value = "`$/saco 25 kilos"
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 300
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
